# The exported sheet used to duplicate the "Название"/"Предмет" columns
# (D/E) alongside the real "Этап"/"Результат" columns (F/G). The export
# codestyle fix drops the redundant duplicate columns entirely, so delete
# worksheet columns D:E outright (shifting the former F:G left into D:E).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1:E1").EntireColumn.Delete()
